$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = 1234588
$ws.Cells.Item(4, 3).Value = 21753
$ws.Cells.Item(4, 4).Value = 199148
$ws.Cells.Item(4, 5).Value = 963386
$ws.Cells.Item(4, 7).Value = 2133
$ws.Cells.Item(4, 8).Value = 72054
$ws.Cells.Item(9, 2).Value = 166706
$ws.Cells.Item(9, 3).Value = 554
$ws.Cells.Item(9, 5).Value = 24613
$ws.Cells.Item(12, 2).Value = 113844
$ws.Cells.Item(12, 3).Value = 5578
$ws.Cells.Item(12, 5).Value = 60181
$ws.Cells.Item(12, 7).Value = 505
$ws.Cells.Item(12, 8).Value = 7848
$ws.Cells.Item(85, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(85, 2).Value = 1464
$ws.Cells.Item(85, 3).Value = 32
$ws.Cells.Item(85, 4).Value = 701
$ws.Cells.Item(85, 5).Value = 745
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 8).Value = 18
$ws.Cells.Item(86, 1).Value = 'Eslovenia'
$ws.Cells.Item(86, 2).Value = 1445
$ws.Cells.Item(86, 3).Value = 6
$ws.Cells.Item(86, 4).Value = 244
$ws.Cells.Item(86, 5).Value = 1103
$ws.Cells.Item(86, 6).Value = 17
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 98
$ws.Cells.Item(93, 2).Value = 1022
$ws.Cells.Item(93, 3).Value = 4
$ws.Cells.Item(93, 4).Value = 482
$ws.Cells.Item(93, 5).Value = 497
$ws.Cells.Item(93, 6).Value = 17
$ws.Cells.Item(108, 1).Value = 'Burkina Faso'
$ws.Cells.Item(108, 2).Value = 688
$ws.Cells.Item(108, 3).Value = 16
$ws.Cells.Item(108, 4).Value = 548
$ws.Cells.Item(108, 5).Value = 92
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 48
$ws.Cells.Item(109, 1).Value = 'Mayotte'
$ws.Cells.Item(109, 2).Value = 686
$ws.Cells.Item(109, 4).Value = 352
$ws.Cells.Item(109, 5).Value = 328
$ws.Cells.Item(109, 6).Value = 6
$ws.Cells.Item(109, 8).Value = 6
$ws.Cells.Item(140, 1).Value = 'Santo Tome y Principe'
$ws.Cells.Item(140, 2).Value = 174
$ws.Cells.Item(140, 3).Value = 151
$ws.Cells.Item(140, 4).Value = 4
$ws.Cells.Item(140, 5).Value = 167
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 3
$ws.Cells.Item(141, 1).Value = 'Republica del Chad'
$ws.Cells.Item(141, 2).Value = 170
$ws.Cells.Item(141, 3).Value = 53
$ws.Cells.Item(141, 4).Value = 43
$ws.Cells.Item(141, 5).Value = 110
$ws.Cells.Item(141, 7).Value = 7
$ws.Cells.Item(141, 8).Value = 17
$ws.Cells.Item(142, 1).Value = 'Liberia'
$ws.Cells.Item(142, 2).Value = 170
$ws.Cells.Item(142, 3).Value = 4
$ws.Cells.Item(142, 4).Value = 58
$ws.Cells.Item(142, 5).Value = 92
$ws.Cells.Item(142, 7).Value = 2
$ws.Cells.Item(142, 8).Value = 20
$ws.Cells.Item(143, 1).Value = 'Birmania'
$ws.Cells.Item(143, 2).Value = 161
$ws.Cells.Item(143, 4).Value = 49
$ws.Cells.Item(143, 5).Value = 106
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 8).Value = 6
$ws.Cells.Item(144, 1).Value = 'Guadalupe'
$ws.Cells.Item(144, 2).Value = 152
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 98
$ws.Cells.Item(144, 5).Value = 42
$ws.Cells.Item(144, 6).Value = 5
$ws.Cells.Item(144, 8).Value = 12
$ws.Cells.Item(145, 1).Value = 'Madagascar'
$ws.Cells.Item(145, 2).Value = 151
$ws.Cells.Item(145, 3).Value = 2
$ws.Cells.Item(145, 4).Value = 101
$ws.Cells.Item(145, 6).Value = 1
$ws.Cells.Item(145, 8).Value = 0
$ws.Cells.Item(146, 1).Value = 'Etiopia'
$ws.Cells.Item(146, 2).Value = 145
$ws.Cells.Item(146, 3).Value = 5
$ws.Cells.Item(146, 4).Value = 91
$ws.Cells.Item(146, 5).Value = 50
$ws.Cells.Item(146, 7).Value = 1
$ws.Cells.Item(146, 8).Value = 4
$ws.Cells.Item(147, 1).Value = 'Gibraltar'
$ws.Cells.Item(147, 2).Value = 144
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 133
$ws.Cells.Item(147, 5).Value = 11
$ws.Cells.Item(147, 8).Value = 0
$ws.Cells.Item(148, 1).Value = 'Zambia'
$ws.Cells.Item(148, 3).Value = 1
$ws.Cells.Item(148, 4).Value = 92
$ws.Cells.Item(148, 5).Value = 43
$ws.Cells.Item(148, 6).Value = 1
$ws.Cells.Item(148, 8).Value = 3
$ws.Cells.Item(149, 1).Value = 'Brunei'
$ws.Cells.Item(149, 2).Value = 138
$ws.Cells.Item(149, 4).Value = 131
$ws.Cells.Item(149, 5).Value = 6
$ws.Cells.Item(150, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(150, 2).Value = 133
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 100
$ws.Cells.Item(150, 5).Value = 32
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 8).Value = 1
$ws.Cells.Item(151, 1).Value = 'Togo'
$ws.Cells.Item(151, 2).Value = 128
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(151, 4).Value = 74
$ws.Cells.Item(151, 5).Value = 45
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 8).Value = 9
$ws.Cells.Item(152, 1).Value = 'Camboya'
$ws.Cells.Item(152, 2).Value = 122
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 120
$ws.Cells.Item(152, 5).Value = 2
$ws.Cells.Item(152, 6).Value = 1
$ws.Cells.Item(152, 8).Value = 0
$ws.Cells.Item(153, 1).Value = 'Suazilandia'
$ws.Cells.Item(153, 2).Value = 119
$ws.Cells.Item(153, 3).Value = 3
$ws.Cells.Item(153, 4).Value = 12
$ws.Cells.Item(153, 5).Value = 106
$ws.Cells.Item(153, 8).Value = 1
$ws.Cells.Item(154, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(154, 2).Value = 116
$ws.Cells.Item(154, 4).Value = 102
$ws.Cells.Item(154, 5).Value = 6
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 8).Value = 8
$ws.Cells.Item(155, 1).Value = 'Bermudas'
$ws.Cells.Item(155, 2).Value = 115
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 54
$ws.Cells.Item(155, 5).Value = 54
$ws.Cells.Item(155, 8).Value = 7
$ws.Cells.Item(156, 1).Value = 'Aruba'
$ws.Cells.Item(156, 2).Value = 101
$ws.Cells.Item(156, 3).Value = 1
$ws.Cells.Item(156, 4).Value = 82
$ws.Cells.Item(156, 5).Value = 17
$ws.Cells.Item(156, 8).Value = 2
$ws.Cells.Item(157, 1).Value = 'Haiti'
$ws.Cells.Item(157, 2).Value = 100
$ws.Cells.Item(157, 4).Value = 10
$ws.Cells.Item(157, 5).Value = 79
$ws.Cells.Item(157, 8).Value = 11
$ws.Cells.Item(158, 1).Value = 'Uganda'
$ws.Cells.Item(158, 2).Value = 97
$ws.Cells.Item(158, 4).Value = 55
$ws.Cells.Item(158, 5).Value = 42
$ws.Cells.Item(158, 8).Value = 0
$ws.Cells.Item(159, 1).Value = 'Benin'
$ws.Cells.Item(159, 2).Value = 96
$ws.Cells.Item(159, 4).Value = 50
$ws.Cells.Item(159, 5).Value = 44
$ws.Cells.Item(159, 8).Value = 2
$ws.Cells.Item(160, 1).Value = 'Monaco'
$ws.Cells.Item(160, 2).Value = 95
$ws.Cells.Item(160, 4).Value = 81
$ws.Cells.Item(160, 5).Value = 10
$ws.Cells.Item(160, 6).Value = 1
$ws.Cells.Item(160, 8).Value = 4
$ws.Cells.Item(161, 1).Value = 'Guyana'
$ws.Cells.Item(161, 2).Value = 92
$ws.Cells.Item(161, 4).Value = 27
$ws.Cells.Item(161, 5).Value = 56
$ws.Cells.Item(161, 6).Value = 3
$ws.Cells.Item(161, 8).Value = 9
$ws.Cells.Item(162, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(162, 2).Value = 85
$ws.Cells.Item(162, 4).Value = 10
$ws.Cells.Item(162, 5).Value = 75
$ws.Cells.Item(162, 8).Value = 0
$ws.Cells.Item(163, 1).Value = 'Bahamas'
$ws.Cells.Item(163, 2).Value = 83
$ws.Cells.Item(163, 4).Value = 25
$ws.Cells.Item(163, 5).Value = 47
$ws.Cells.Item(163, 6).Value = 1
$ws.Cells.Item(163, 8).Value = 11
$ws.Cells.Item(164, 1).Value = 'Nepal'
$ws.Cells.Item(164, 3).Value = 7
$ws.Cells.Item(164, 4).Value = 16
$ws.Cells.Item(164, 5).Value = 66
$ws.Cells.Item(164, 8).Value = 0
$ws.Cells.Item(165, 1).Value = 'Barbados'
$ws.Cells.Item(165, 4).Value = 47
$ws.Cells.Item(165, 5).Value = 28
$ws.Cells.Item(165, 6).Value = 4
$ws.Cells.Item(165, 8).Value = 7
$ws.Cells.Item(166, 1).Value = 'Liechtenstein'
$ws.Cells.Item(166, 2).Value = 82
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 8).Value = 1
$ws.Cells.Item(167, 1).Value = 'Mozambique'
$ws.Cells.Item(167, 2).Value = 81
$ws.Cells.Item(167, 3).Value = 1
$ws.Cells.Item(168, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(168, 2).Value = 76
$ws.Cells.Item(168, 4).Value = 44
$ws.Cells.Item(168, 5).Value = 18
$ws.Cells.Item(168, 6).Value = 7
$ws.Cells.Item(168, 7).Value = 1
$ws.Cells.Item(168, 8).Value = 14
$ws.Cells.Item(169, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(169, 2).Value = 75
$ws.Cells.Item(169, 4).Value = 14
$ws.Cells.Item(169, 5).Value = 60
$ws.Cells.Item(169, 6).Value = 3
$ws.Cells.Item(169, 8).Value = 1
$ws.Cells.Item(170, 1).Value = 'Libia'
$ws.Cells.Item(170, 2).Value = 63
$ws.Cells.Item(170, 4).Value = 23
$ws.Cells.Item(170, 5).Value = 37
$ws.Cells.Item(170, 8).Value = 3
$ws.Cells.Item(171, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(171, 2).Value = 58
$ws.Cells.Item(171, 4).Value = 53
$ws.Cells.Item(171, 5).Value = 5
$ws.Cells.Item(171, 6).Value = 1
$ws.Cells.Item(172, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(172, 2).Value = 52
$ws.Cells.Item(172, 3).Value = 6
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 52
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(173, 1).Value = 'Macao'
$ws.Cells.Item(173, 2).Value = 45
$ws.Cells.Item(173, 4).Value = 39
$ws.Cells.Item(173, 5).Value = 6
$ws.Cells.Item(173, 6).Value = 1
$ws.Cells.Item(174, 1).Value = 'Siria'
$ws.Cells.Item(174, 2).Value = 44
$ws.Cells.Item(174, 4).Value = 27
$ws.Cells.Item(174, 5).Value = 14
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 8).Value = 0
$ws.Cells.Item(175, 1).Value = 'Malaui'
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 9
$ws.Cells.Item(175, 5).Value = 29
$ws.Cells.Item(175, 6).Value = 1
$ws.Cells.Item(175, 8).Value = 3
$ws.Cells.Item(176, 1).Value = 'Mongolia'
$ws.Cells.Item(176, 2).Value = 41
$ws.Cells.Item(176, 3).Value = 1
$ws.Cells.Item(176, 4).Value = 13
$ws.Cells.Item(176, 5).Value = 28
$ws.Cells.Item(176, 8).Value = 0
$ws.Cells.Item(177, 1).Value = 'Puerto Rico'
$ws.Cells.Item(177, 4).Value = 1
$ws.Cells.Item(177, 5).Value = 36
$ws.Cells.Item(177, 8).Value = 2
$ws.Cells.Item(178, 1).Value = 'Eritrea'
$ws.Cells.Item(178, 2).Value = 39
$ws.Cells.Item(178, 4).Value = 30
$ws.Cells.Item(178, 5).Value = 9
$ws.Cells.Item(178, 8).Value = 0
$ws.Cells.Item(179, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(179, 2).Value = 38
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 29
$ws.Cells.Item(179, 5).Value = 6
$ws.Cells.Item(179, 6).Value = 1
$ws.Cells.Item(179, 8).Value = 3
$ws.Cells.Item(180, 1).Value = 'Angola'
$ws.Cells.Item(180, 2).Value = 36
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(180, 4).Value = 11
$ws.Cells.Item(180, 5).Value = 23
$ws.Cells.Item(180, 8).Value = 2
$ws.Cells.Item(181, 1).Value = 'Zimbabue'
$ws.Cells.Item(181, 2).Value = 34
$ws.Cells.Item(181, 4).Value = 5
$ws.Cells.Item(181, 5).Value = 25
$ws.Cells.Item(181, 8).Value = 4
$ws.Cells.Item(182, 1).Value = 'Guam'
$ws.Cells.Item(182, 2).Value = 32
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 31
$ws.Cells.Item(182, 8).Value = 1
$ws.Cells.Item(183, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(183, 2).Value = 25
$ws.Cells.Item(183, 4).Value = 16
$ws.Cells.Item(183, 5).Value = 6
$ws.Cells.Item(183, 6).Value = 1
$ws.Cells.Item(183, 8).Value = 3
$ws.Cells.Item(184, 1).Value = 'Timor Oriental'
$ws.Cells.Item(184, 2).Value = 24
$ws.Cells.Item(184, 4).Value = 20
$ws.Cells.Item(184, 5).Value = 4
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(186, 2).Value = 22
$ws.Cells.Item(186, 3).Value = 10
$ws.Cells.Item(186, 5).Value = 18
$ws.Cells.Item(197, 1).Value = 'Nicaragua'
$ws.Cells.Item(197, 3).Value = 1
$ws.Cells.Item(197, 4).Value = 7
$ws.Cells.Item(197, 5).Value = 4
$ws.Cells.Item(197, 8).Value = 5
$ws.Cells.Item(198, 1).Value = 'Curazao'
$ws.Cells.Item(198, 4).Value = 13
$ws.Cells.Item(198, 8).Value = 1
$ws.Cells.Item(199, 1).Value = 'Dominica'
$ws.Cells.Item(199, 2).Value = 16
$ws.Cells.Item(199, 4).Value = 14
$ws.Cells.Item(199, 5).Value = 2
$ws.Cells.Item(200, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(200, 4).Value = 8
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(201, 1).Value = 'Burundi'
$ws.Cells.Item(201, 5).Value = 7
$ws.Cells.Item(201, 8).Value = 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 6 de Mayo de 2020 a las 00:03'
